# Add new columns I (I0) and J (IF) to the sheet, mirroring existing
# header/style conventions used by the other columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells: same bold/centered/bordered style as the other headers (s="1")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data values for rows 2-41 (I and J columns)
$iValues = @(9,6,7,8,7,8,5,7,8,6,6,1,7,6,7,7,6,8,7,4,9,7,6,6,9,7,6,8,7,5,9,7,7,6,7,8,7,5,7,5)
$jValues = @(9,6,7,8,7,8,5,7,8,6,6,1,7,6,7,7,6,8,7,4,9,7,6,7,9,7,7,8,7,5,9,7,7,6,7,8,7,5,7,5)

for ($k = 0; $k -lt $iValues.Length; $k++) {
    $row = $k + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$k]
    $ws.Cells.Item($row, 10).Value = $jValues[$k]
}
